# regen sval data to filter save games
# Update the raw stat columns (B:E) for each game row; column G ("sum")
# is re-derived as SUM(B:E) for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    3  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    4  = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 13.86384647080068)
    5  = @(3.272327238179451, 9.983522426115931, 0.7210945179870265, 13.86384647080068)
    6  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    7  = @(0.2881169905109251, 1.626987699542094, 18.71679738969934, 13.86384647080068)
    8  = @(3.272327238179451, 0.3048912486333797, 3.223369029078222, 0.5333859586016987)
    9  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    10 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    11 = @(0.1169995834814548, 9.983522426115931, 0.7210945179870265, 13.86384647080068)
    12 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 13.86384647080068)
    13 = @(1.445647641019636, 1.626987699542094, 18.71679738969934, 0.5333859586016987)
    14 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    15 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068)
    16 = @(0.1169995834814548, 1.626987699542094, 18.71679738969934, 13.86384647080068)
    17 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $sum = $vals[0] + $vals[1] + $vals[2] + $vals[3]
    $ws.Cells.Item($row, 7).Value = $sum
}
